$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear early-year forecast values that are no longer computed
$ws.Range("C3:C6").ClearContents()
$ws.Range("E2:E6").ClearContents()

# Update the y_0_forecast (C) and y_1_forecast (E) values for rows 7-19
$ws.Range("C7").Value = -0.06726240733891942
$ws.Range("C8").Value = -0.3584575688954428
$ws.Range("C9").Value = -0.009688851112665819
$ws.Range("C10").Value = -0.1132037832954791
$ws.Range("C11").Value = -0.2638577853126156
$ws.Range("C12").Value = 0.04166709579394023
$ws.Range("C13").Value = -0.7671134292608239
$ws.Range("C14").Value = -0.2043373675692961
$ws.Range("C15").Value = -0.8080927309597863
$ws.Range("C16").Value = -1.697148566375528
$ws.Range("C17").Value = -0.6247846736575413
$ws.Range("C18").Value = -0.09040308684795662
$ws.Range("C19").Value = 0.3813481955213138

$ws.Range("E7").Value = -0.4774698422615242
$ws.Range("E8").Value = -0.4112436562971
$ws.Range("E9").Value = -0.3537462851234685
$ws.Range("E10").Value = -0.3204027102583273
$ws.Range("E11").Value = -0.3143428957755656
$ws.Range("E12").Value = -0.2532926704812977
$ws.Range("E13").Value = -0.3480983700859808
$ws.Range("E14").Value = -0.1872494743064723
$ws.Range("E15").Value = -0.1644030883838465
$ws.Range("E16").Value = -0.2339923140600275
$ws.Range("E17").Value = -0.3554771869619944
$ws.Range("E18").Value = -0.3159097170635006
$ws.Range("E19").Value = -0.2523778956734835
